# Bournemouth_stats.xlsx edit
#  1) Rename the stat sheets to their "friendly" (spaced) names.
#  2) Bump the "days" component of every player's Age column (col E,
#     format "YY-DDD") by one day, on every stats sheet (everything
#     except the "Matches" sheet).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $newName
}

# --- 2. Bump the Age "days" counter on every stats sheet --------------
foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Matches") {
        continue
    }

    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $val = $cell.Text

        if ($val -match "^(\d+)-(\d+)$") {
            $years = $matches[1]
            $days = [int]$matches[2]
            $days = $days + 1
            $newVal = $years + "-" + $days.ToString("D3")
            $cell.Value = $newVal
        }
    }
}

Write-Output "done"
